$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F5").Value = 357
$ws1.Range("F6").Value = 150
$ws1.Range("F7").Value = 152
$ws1.Range("F8").Value = 747
$ws1.Range("F9").Value = 4104
$ws1.Range("F11").Value = 52
$ws1.Range("F12").Value = 165
$ws1.Range("F13").Value = 95
$ws1.Range("F14").Value = 5854
$ws1.Range("F15").Value = 458
$ws1.Range("F16").Value = 2273
$ws1.Range("F17").Value = 256
$ws1.Range("F18").Value = 155
$ws1.Range("F19").Value = 439
$ws1.Range("F20").Value = 8821
$ws1.Range("F21").Value = 39
$ws1.Range("F22").Value = 1500
$ws1.Range("F23").Value = 192
$ws1.Range("F24").Value = 2275
$ws1.Range("F25").Value = 2353
$ws1.Range("F26").Value = 1371
$ws1.Range("F27").Value = 218
$ws1.Range("F28").Value = 1917
$ws1.Range("F29").Value = 35
$ws1.Range("F31").Value = 317
$ws1.Range("F33").Value = 30
$ws1.Range("F34").Value = 277
$ws1.Range("F36").Value = 13
$ws1.Range("F37").Value = 27
$ws1.Range("F39").Value = 1205
$ws1.Range("F40").Value = 63
$ws1.Range("F41").Value = 86
$ws1.Range("F42").Value = 222
$ws1.Range("F43").Value = 1475
$ws1.Range("F44").Value = 2375
$ws1.Range("F46").Value = 899
$ws1.Range("F47").Value = 278
$ws2.Range("F2").Value = 5
$ws2.Range("F5").Value = 166
$ws2.Range("F12").Value = 144
$ws3.Range("F2").Value = 671
$ws3.Range("F3").Value = 865
$ws3.Range("F4").Value = 91
$ws4.Range("F4").Value = 671
$ws4.Range("F5").Value = 865
$ws4.Range("F6").Value = 91
$ws4.Range("F7").Value = 357
$ws4.Range("F8").Value = 5
$ws4.Range("F9").Value = 150
$ws4.Range("F11").Value = 152
$ws4.Range("F12").Value = 747
$ws4.Range("F13").Value = 4104
$ws4.Range("F14").Value = 4104
$ws4.Range("F15").Value = 52
$ws4.Range("F16").Value = 165
$ws4.Range("F18").Value = 5854
$ws4.Range("F19").Value = 458
$ws4.Range("F20").Value = 2273
$ws4.Range("F22").Value = 155
$ws4.Range("F23").Value = 439
$ws4.Range("F24").Value = 8821
$ws4.Range("F25").Value = 144
$ws4.Range("F26").Value = 39
$ws4.Range("F27").Value = 1501
$ws4.Range("F28").Value = 2275
$ws4.Range("F29").Value = 2353
$ws4.Range("F30").Value = 1371
$ws4.Range("F31").Value = 218
$ws4.Range("F32").Value = 1917
$ws4.Range("F33").Value = 35
$ws4.Range("F35").Value = 317
$ws4.Range("F36").Value = 277
$ws4.Range("F38").Value = 27
$ws4.Range("F40").Value = 63
$ws4.Range("F41").Value = 222
$ws4.Range("F42").Value = 1475
$ws4.Range("F43").Value = 2375
$ws4.Range("F44").Value = 899
$ws4.Range("F46").Value = 278
